# [POD-534] Rewrote tests broken by new binnable/boxable rules
#
# The fixture's "Bin barcode" (A) and "Bin identifier" (B) sample values are
# refreshed, and a new "Box barcode" (C) sample value is populated for both
# data rows. The active selection moves to A2, column E is widened, and the
# recalc id is bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2, 5) {
    $ws.Cells.Item($r, 1).Value = 40000000000028   # A: Bin barcode
    $ws.Cells.Item($r, 2).Value = "test-bin"        # B: Bin identifier
    $ws.Cells.Item($r, 3).Value = 40000000000036   # C: Box barcode
}

$ws.Columns.Item(5).ColumnWidth = 13.166666666666666   # xlsx col width 14

$ws.Range("A2").Select()

$wb.CalcId = 152511
